$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (I1, J1) - copy formatting from H1 (bold, bordered, centered header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows (plain numeric, no special style)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4

$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1

$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 6
